# Insert a new weekly price record at row 22 ("Hortaliza, Terminal La
# Palmera de La Serena - Poroto granado"). Excel's native row-insert
# shifts the existing rows 22-88 down to 23-89 (and grows the used range
# to A1:R89), after which we populate the newly-inserted row with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(22).Insert()

$ws.Cells.Item(22, 1).Value  = 8
$ws.Cells.Item(22, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(22, 3).Value  = "Coquimbo"
$ws.Cells.Item(22, 4).Value  = 44659
$ws.Cells.Item(22, 5).Value  = 4
$ws.Cells.Item(22, 6).Value  = 100112030
$ws.Cells.Item(22, 7).Value  = "Poroto granado"
$ws.Cells.Item(22, 8).Value  = "Sin especificar"
$ws.Cells.Item(22, 9).Value  = "Primera"
$ws.Cells.Item(22, 10).Value = 400
$ws.Cells.Item(22, 11).Value = 27000
$ws.Cells.Item(22, 12).Value = 28000
$ws.Cells.Item(22, 13).Value = 27500
$ws.Cells.Item(22, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(22, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(22, 16).Value = 1100
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"
